# Add Test Scenario IDs (TS ID) for the User Registration rows (FR_REG_01..FR_REG_10)
# on the RTM_ALL sheet, column D, rows 11-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM_ALL")

# Entry order mirrors how the values were actually typed into the sheet so the
# resulting shared-string table matches the authored edit: rows 11-14 reuse
# existing strings, then rows 18-20 are filled in, followed by rows 15-17.
$ws.Range("D11").Value = "TS_REG_01"
$ws.Range("D12").Value = "TS_REG_02"
$ws.Range("D13").Value = "TS_REG_03"
$ws.Range("D14").Value = "TS_REG_04"
$ws.Range("D18").Value = "TS_REG_08"
$ws.Range("D19").Value = "TS_REG_09"
$ws.Range("D20").Value = "TS_REG_10"
$ws.Range("D15").Value = "TS_REG_05"
$ws.Range("D16").Value = "TS_REG_06"
$ws.Range("D17").Value = "TS_REG_07"

# Make this the active sheet and set the selection to match the authored edit.
$ws.Activate()
$ws.Range("E12").Select()
